$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing Ruilin annotation row (row 2) into a new row 3,
# which keeps the original cell types (e.g. the text "3" for politeness_score).
$ws.Range("A2:H2").Copy()
$ws.Range("A3").PasteSpecial()

# Row 3 gets the second annotation's differing fields (issue_type, id,
# source_file, text) while Annotator/politeness_score/polite_expressions/
# sentence_purpose stay the same as row 2.
$ws.Range("E3").Value = "WRI"
$ws.Range("F3").Value = "c8048836-24fe-4e27-95aa-c7cfb58ac155"
$ws.Range("G3").Value = "rkc_hGb0Z_annotated.xlsx"
$ws.Range("H3").Value = "The structure of the global policies used in the experiments should be mentioned somewhere."

# Row 2's politeness_score is normalized to a true number.
$ws.Range("B2").Value = 3
